{"js": "// Edit: add \"Double hashing\" section, split \"quadratic probing\" run for spellcheck,\n// and relocate the _GoBack bookmark from the final (now replaced) empty paragraph\n// into the \"Hashfunktionen ... \u00e4knar ...\" paragraph, per the commit\n// \"lagt tillstycke om double hashing\" (added a paragraph about double hashing).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Change 1: split \" quadratic probing \" run so \"quadratic\" is flagged spellStart/spellEnd ---\nconst quadraticParaXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>quadratic</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> probing </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nparagraphs.items[5].insertOoxml(quadraticParaXml, Word.InsertLocation.replace);\n\n// --- Change 2: move the _GoBack bookmark into the \"Hashfunktionen ... \u00e4knar ...\" paragraph ---\nconst hashfunktionParaXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>Hashkoden</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> skapas med en </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashfunktion</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\">. </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>Hashfunktionen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:t>r</w:t></w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n<w:r><w:t xml:space=\"preserve\">\u00e4knar om elementets data till en kod. </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>Hashfunktionen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> ska alltid ge samma kod f\u00f6r samma objekt. En mindre effektiv </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashkod</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> g\u00f6r ofta att flera element hamnar p\u00e5 samma plats, medan en mer effektiv </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashfunktion</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> g\u00f6r att elementen sprids ut j\u00e4mt \u00f6ver hela </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>arrayen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\">, och d\u00e4rmed minimerar risken f\u00f6r krockar. </w:t></w:r>\n</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nparagraphs.items[10].insertOoxml(hashfunktionParaXml, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// --- Change 3: replace the final (bookmark-only) paragraph with the new \"Double hashing\" text ---\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst doubleHashingParaXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>\n<w:r><w:t xml:space=\"preserve\">Double </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashing</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> \u00e4r ytterligare en metod f\u00f6r att hantera kollisioner</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">. Double </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashing</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> inneb\u00e4r att en </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashkod</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> tas fram utifr\u00e5n elementet som ska placeras. </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">Om en kollision uppst\u00e5r p\u00e5 denna plats, r\u00e4knas en ny </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashkod</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> ut, </w:t></w:r>\n<w:r><w:t>f\u00f6r det interval</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">l som ska hoppas fram i </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>arrayen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t>.</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> P\u00e5 denna nya plats kan sedan elementet stoppas om, om den \u00e4r ledig. </w:t></w:r>\n</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nlastParagraph.insertOoxml(doubleHashingParaXml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Edit: add \"Double hashing\" section, split \"quadratic probing\" run for spellcheck,\n# and relocate the _GoBack bookmark from the final (now replaced) empty paragraph\n# into the \"Hashfunktionen ... \u00e4knar ...\" paragraph, per the commit\n# \"lagt tillstycke om double hashing\" (added a paragraph about double hashing).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: split \" quadratic probing \" run so \"quadratic\" is flagged spellStart/spellEnd ---\n$quadraticParaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>quadratic</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> probing </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs.Item(6).Range.InsertXML($quadraticParaXml) | Out-Null\n\n# --- Change 2: move the _GoBack bookmark into the \"Hashfunktionen ... \u00e4knar ...\" paragraph ---\n$hashfunktionParaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>Hashkoden</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> skapas med en </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashfunktion</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\">. </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>Hashfunktionen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n<w:r><w:t>r</w:t></w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n<w:r><w:t xml:space=\"preserve\">\u00e4knar om elementets data till en kod. </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>Hashfunktionen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> ska alltid ge samma kod f\u00f6r samma objekt. En mindre effektiv </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashkod</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> g\u00f6r ofta att flera element hamnar p\u00e5 samma plats, medan en mer effektiv </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashfunktion</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> g\u00f6r att elementen sprids ut j\u00e4mt \u00f6ver hela </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>arrayen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\">, och d\u00e4rmed minimerar risken f\u00f6r krockar. </w:t></w:r>\n</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs.Item(11).Range.InsertXML($hashfunktionParaXml) | Out-Null\n\n# --- Change 3: replace the final (bookmark-only) paragraph with the new \"Double hashing\" text ---\n$lastIndex = $d.Paragraphs.Count\n$doubleHashingParaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>\n<w:r><w:t xml:space=\"preserve\">Double </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashing</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> \u00e4r ytterligare en metod f\u00f6r att hantera kollisioner</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">. Double </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashing</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> inneb\u00e4r att en </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashkod</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> tas fram utifr\u00e5n elementet som ska placeras. </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">Om en kollision uppst\u00e5r p\u00e5 denna plats, r\u00e4knas en ny </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>hashkod</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t xml:space=\"preserve\"> ut, </w:t></w:r>\n<w:r><w:t>f\u00f6r det interval</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">l som ska hoppas fram i </w:t></w:r>\n<w:proofErr w:type=\"spellStart\"/>\n<w:r><w:t>arrayen</w:t></w:r>\n<w:proofErr w:type=\"spellEnd\"/>\n<w:r><w:t>.</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> P\u00e5 denna nya plats kan sedan elementet stoppas om, om den \u00e4r ledig. </w:t></w:r>\n</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$d.Paragraphs.Item($lastIndex).Range.InsertXML($doubleHashingParaXml) | Out-Null\n\nWrite-Output \"done\"\n"}
